$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "MIN"
$ws.Range("B11").Formula = "=MIN(B4:B8)"

$ws.Range("A12").Value = "MAX"
$ws.Range("B12").Formula = "=MAX(B4:B8)"

$ws.Range("A13").Value = "AVERAGE"

$ws.Range("G11").Select()
